$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Mfge8"
$ws.Range("C2").Value = "Pdgfrb"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 21.66138833333333
$ws.Range("H2").Value = 64.98416499999999
$ws.Range("I2").Value = 0.2252765553546639
$ws.Range("J2").Value = 0.2252765553546639
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 7.259429666666667
$ws.Range("N2").Value = 21.778289
$ws.Range("O2").Value = 0.05296410708422108
$ws.Range("P2").Value = 0.05296410708422109
$ws.Range("Q2").Value = 157.2493250881872
$ws.Range("R2").Value = 1415.243925793685
$ws.Range("S2").Value = 0.01193157160136888
$ws.Range("T2").Value = 0.01193157160136888

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Mfge8"
$ws.Range("C3").Value = "Pdgfrb"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 21.66138833333333
$ws.Range("H3").Value = 64.98416499999999
$ws.Range("I3").Value = 0.2252765553546639
$ws.Range("J3").Value = 0.2252765553546639
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 60.13240533333334
$ws.Range("N3").Value = 180.397216
$ws.Range("O3").Value = 0.438720299189682
$ws.Range("P3").Value = 0.4387202991896821
$ws.Range("Q3").Value = 1302.551383342738
$ws.Range("R3").Value = 11722.96245008464
$ws.Range("S3").Value = 0.0988333977656191
$ws.Range("T3").Value = 0.09883339776561913

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Mfge8"
$ws.Range("C4").Value = "Pdgfrb"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 21.66138833333333
$ws.Range("H4").Value = 64.98416499999999
$ws.Range("I4").Value = 0.2252765553546639
$ws.Range("J4").Value = 0.2252765553546639
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 69.67135866666666
$ws.Range("N4").Value = 209.014076
$ws.Range("O4").Value = 0.5083155937260968
$ws.Range("P4").Value = 0.5083155937260969
$ws.Range("Q4").Value = 1509.178355789615
$ws.Range("R4").Value = 13582.60520210654
$ws.Range("S4").Value = 0.1145115859876759
$ws.Range("T4").Value = 0.1145115859876759

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Mfge8"
$ws.Range("C5").Value = "Pdgfrb"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 24.68088566666666
$ws.Range("H5").Value = 74.04265699999999
$ws.Range("I5").Value = 0.2566790650963491
$ws.Range("J5").Value = 0.2566790650963491
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 7.259429666666667
$ws.Range("N5").Value = 21.778289
$ws.Range("O5").Value = 0.05296410708422108
$ws.Range("P5").Value = 0.05296410708422109
$ws.Range("Q5").Value = 179.1691536082081
$ws.Range("R5").Value = 1612.522382473873
$ws.Range("S5").Value = 0.01359477749004079
$ws.Range("T5").Value = 0.01359477749004079

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Mfge8"
$ws.Range("C6").Value = "Pdgfrb"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 24.68088566666666
$ws.Range("H6").Value = 74.04265699999999
$ws.Range("I6").Value = 0.2566790650963491
$ws.Range("J6").Value = 0.2566790650963491
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 60.13240533333334
$ws.Range("N6").Value = 180.397216
$ws.Range("O6").Value = 0.438720299189682
$ws.Range("P6").Value = 0.4387202991896821
$ws.Range("Q6").Value = 1484.121020893657
$ws.Range("R6").Value = 13357.08918804291
$ws.Range("S6").Value = 0.1126103162347981
$ws.Range("T6").Value = 0.1126103162347982

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Mfge8"
$ws.Range("C7").Value = "Pdgfrb"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 24.68088566666666
$ws.Range("H7").Value = 74.04265699999999
$ws.Range("I7").Value = 0.2566790650963491
$ws.Range("J7").Value = 0.2566790650963491
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 69.67135866666666
$ws.Range("N7").Value = 209.014076
$ws.Range("O7").Value = 0.5083155937260968
$ws.Range("P7").Value = 0.5083155937260969
$ws.Range("Q7").Value = 1719.550837493326
$ws.Range("R7").Value = 15475.95753743993
$ws.Range("S7").Value = 0.1304739713715101
$ws.Range("T7").Value = 0.1304739713715102

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Mfge8"
$ws.Range("C8").Value = "Pdgfrb"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 49.812376
$ws.Range("H8").Value = 149.437128
$ws.Range("I8").Value = 0.518044379548987
$ws.Range("J8").Value = 0.5180443795489871
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 7.259429666666667
$ws.Range("N8").Value = 21.778289
$ws.Range("O8").Value = 0.05296410708422108
$ws.Range("P8").Value = 0.05296410708422109
$ws.Range("Q8").Value = 361.6094401015547
$ws.Range("R8").Value = 3254.484960913992
$ws.Range("S8").Value = 0.02743775799281142
$ws.Range("T8").Value = 0.02743775799281143

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Mfge8"
$ws.Range("C9").Value = "Pdgfrb"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 49.812376
$ws.Range("H9").Value = 149.437128
$ws.Range("I9").Value = 0.518044379548987
$ws.Range("J9").Value = 0.5180443795489871
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 60.13240533333334
$ws.Range("N9").Value = 180.397216
$ws.Range("O9").Value = 0.438720299189682
$ws.Range("P9").Value = 0.4387202991896821
$ws.Range("Q9").Value = 2995.337984248406
$ws.Range("R9").Value = 26958.04185823565
$ws.Range("S9").Value = 0.2272765851892647
$ws.Range("T9").Value = 0.2272765851892648

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Mfge8"
$ws.Range("C10").Value = "Pdgfrb"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 49.812376
$ws.Range("H10").Value = 149.437128
$ws.Range("I10").Value = 0.518044379548987
$ws.Range("J10").Value = 0.5180443795489871
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 69.67135866666666
$ws.Range("N10").Value = 209.014076
$ws.Range("O10").Value = 0.5083155937260968
$ws.Range("P10").Value = 0.5083155937260969
$ws.Range("Q10").Value = 3470.495914334858
$ws.Range("R10").Value = 31234.46322901373
$ws.Range("S10").Value = 0.2633300363669108
$ws.Range("T10").Value = 0.2633300363669109

Write-Output "done"